$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SOFT130015"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 17302010049
$ws.Range("D2").Value = "A"

$ws.Range("A3").Value = "BIOL110007"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 17302010049
$ws.Range("D3").Value = "B+"

$ws.Range("A4").Value = "SOFT130049"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 17302010049
$ws.Range("D4").Value = "C"

$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

$ws.Range("D9").Select() | Out-Null
